$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40 numeric data
$ws.Range("A40").Value = 45888
$ws.Range("A40").NumberFormat = $ws.Range("A39").NumberFormat
$ws.Range("B40").Value = 6411.3701171875
$ws.Range("C40").Value = 6456.47998046875
$ws.Range("D40").Value = 6400.22021484375
$ws.Range("E40").Value = 6446.240234375
$ws.Range("F40").Value = 4377440000
$ws.Range("G40").Value = -0.0058581031187567
$ws.Range("H40").Value = -1

# Row 40 headline text columns
$ws.Range("I40").Value = "poriferous llc achieves coveted eu mdr regulatory certification for porous polyethylene implant warehouse data reveals how shipper managed inventory through tariff nvidia earnings face high bar on datacenter demand china deal newsmax to pay million to settle dominion suit over election fraud claim bmws next big suv is gunning for the top dog but will it go electric jules witcover political reporter and columnist dy at samsung overtakes apple u smartphone market update hint at major shift brent wti edge higher a zelensky arrives in washington for talk with trump top trending coin this week oz presale token gain early momentum best altcoins to buy for september remittix top list over ada pepe and bonk cold wallet offer referral bonus bonk burn t supply ripple eye crypto project set to explode top crypto executive seen accumulating million of token of new meme coin challenging shiba inu shib in walmart is selling a queen mattress topper for only thats soft and cool to the touch best electronic signature apps in judge denies bail to north andover offduty officer shot in home defense file another appeal whale bet big on this new payment token to replace xrp and stellar by carvanas breakout fuel interest in consumer ecommerce etf why are colorado rabbit growing tentacle and horn maluma sign with wme in all area the sky is the limit dow jones today stock tick lower a major index hover near record high after two consecutive week of gain investopedia fake labubu doll pose threat to kid safety watchdog say cbs news court of federal claim side with small business in landmark taa v baa protest employer guidance for new tax deduction for tip and overtime pay terra balcanica intersects gold mineralization over m and start drill program at brezani target in bosnia ap business summarybrief at pm edt conservative network newsmax agrees to pay m in defamation case over bogus election claim invesco qqq qqq share sold by lee financial co bank stock to keep an eye on august th terawulf stock jump more than a google boost stake in data center operator top bse gainer to watch in uncy investor have opportunity to lead unicycive therapeutic inc security fraud lawsuit with the schall law firm cornerstone select advisor llc reduces stock position in pfizer inc pfe southeast toyota finance ranked highest in overall dealer satisfaction by jd power forget ai robotic process automation is healthcare hottest technology warehouse teamster at unfi win lucrative first contract promising retail stock to follow now august th alphabet inc goog share sold by wellington management group llp promising lithium stock to keep an eye on august th insurance stock to watch today august th hydrogen stock to consider august th pwc kunai acquisition to bolster tech enablement for bank cfo the bull case for nvidia just got stronger according to these analyst this curved samsung screen is at one of it lowest price today soho house to be taken private in billion deal the new york time expect fiery outcome from fed jackson hole conference this week financial advisor network inc take position in cocacola company the ko what analyst rating have to say about the campbell mesa laboratory executes previouslyannounced strategic financing plan ibm option trading a deep dive into market sentiment forecasting the future analyst projection for cloudflare a closer look at walmarts option market dynamic flex meal deal return for faculty staff and graduate and professional student paypal holdingss option frenzy what you need to know a glimpse into the expert outlook on jm smucker through analyst assessing microsoft insight from financial analyst this ethereum trader turned into million in month the analyst verdict fedex in the eye of expert ai take apple edge foxconns pivot could spark etf shift where prudential financial stand with analyst mesa laboratory executes previouslyannounced strategic financing plan wired roundup why gpt flopped wired is the market bullish or bearish on royal gold looking into willis tower watson recent short interest a look into pge inc price over earnings how is the market feeling about lumentum holding ozempic creator novo nordisk will sell weight loss drug for cheaper through this new partnership novo nordisk offer ozempic at per month to eligible u cashpaying customer reuters texas may have to shut down data center to protect it energy grid can blockdag reach analyst review of the m presale in dogecoin news today dogecoin price decline but technical indicator suggest potential rebound applied digitals stock apld jump on plan for billion ai campus tipranks worker are job hugging in a stagnant labor market but growing resentment mean they could bail a soon a the next great resignation come yahoocom texas attorney general accuses meta characterai of misleading kid with mental health claim techcrunch worker are job hugging in a stagnant labor market but growing resentment mean they could bail a soon a the next great resignation come fortune conservative network newsmax agrees to pay m in defamation case over bogus election claim skechers named title sponsor of world champion cup shriners childrens announced a official charity partner feather sound country club in clearwater florida to host tournament whats going on with propanc share monday mastercard and worldpay team to streamline uae money movement mcdonalds corporation mcd share sold by community bank trust waco texas granite harbor advisor inc boost holding in mcdonalds corporation mcd walmart ha everything you need to head back to school from backpack to laptop madd join force with law enforcement to prevent impaired driving ahead of labor day with saturation saturday initiative this weekend ibd live qa stock list for monday aug spotlight on fed what powell bowman and waller will say about rate dogecoin rally spark buzz but trader say this sub coin ha more upside fairscale capital llc purchase share of crowdstrike crwd lee financial co sell share of danaher corporation dhr berger montague pc investigates security claim against kindercare learning company inc nyse klc amazon is selling an bluetooth speaker for and it perfect for any setting inventhelp inventor develops new training aid for basketball player chk top data science recruiter in find ai analytics and ml expert gsg advisor llc purchase share of invesco qqq qqq gsg advisor llc ha stock holding in mastercard incorporated ma weekly bankruptcy alert august for the week ending august jackson hole loom a powell speech test september ratecut hope newsmax to pay million to settle dominion defamation lawsuit intel stock drop on report of equity stake via the chip act barrons liquid edge ai platform leap expands support to laptop with bestinclass performance on amd ryzentm and ryzen aitm processor ap business summarybrief at pm edt ap business summarybrief at pm edt texas democrat are back what now froma harrop are parent who dont vaccinate their child guilty of abuse handwritten proof of holocaust theft should compel congress to act ap news summary at pm edt ap news summary at pm edt state and dc sue doj to stop immigration requirement on victim fund trump zelenskyy meeting not end of the road for u support in securing a peace deal editorial roundup united state eric adam roll out nypd expansion say trump dc crackdown not a factor ap trending summarybrief at pm edt texas lawmaker in illinois heading back home following twoweek walkout over remap fight nfl team refuse to back down after facing backlash for hiring male cheerleader stephen miller say trump dc crackdown will also target graffiti ap news summary at pm edt conservative network newsmax agrees to pay m in defamation case over bogus election claim mark hamill say his wife clever thought stopped him from fleeing u after trump won trump zelensky put on the smile at oval office vancouver mayor warns delay on i bridge replacement cost m per day trump zelensky put on the smile at oval office hamas accepts an arab ceasefire proposal on gaza a palestinian death toll pass texas redistricting fight shift a democrat return to austin with a new strategy mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown on putin advice trump launch assault on mailin ballot and voting machine eo seek to expand access to alternative asset investment in retirement plan doj demand sensitive illinois voter registration data mike mcgrew an alternative mindset for patriotic independent commentary ha brian stelter doesnt get that he only making msnbcs rebranding even funnier in his excited thread former attorney general william barr say he didnt see trump implicated in epstein file former attorney general william barr say he didnt see trump implicated in epstein file gavin newsom left fox news speechless with his unexpectedly hilarious woke comeback these state leader might be the last line of defense against trump abuse of power subsidence isnt just an environmental crisis and it can be slowed the latest zelenskyy signal openness to threeway meeting with trump and putin doj wont appeal judge order on concealed carry ban in post office the latest zelenskyy signal openness to threeway meeting with trump and putin letter political treachery versus personal charity what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work bnpl provider ask for extension to nydfs request for information deadline wyoming teen educate youth about human trafficking on social medium congress must correct an injustice facing combatinjured veteran how blessing and disaster shape alabama why is the c galaxy such an iconic aircraft texas attorney general accuses meta characterai of misleading kid with mental health claim texas ag accuses meta characterai of misleading kid with mental health claim european leader join trump zelenskyy at white house conservative network newsmax agrees to pay m in defamation case over bogus election claim what to know about redistricting fight a texas democrat return and california start work what to know about redistricting fight a texas democrat return and california start work twoweek standoff end quorum reached in texas house mississippi becomes fourth state to send national guard troop to dc in expanding federal crackdown letter trumpputin photo op accomplished nothing judge weighs detainee legal right at alligator alcatraz in florida everglades judge weighs detainee legal right at alligator alcatraz in florida everglades roblox is locking down sexual content and access to adult location after lawsuit trump say the u would back european security guarantee for ukraine to help end the war with russia good in a crisis dr nirav shah could be a promising candidate for anything steve collins wear a suit bro zelenskyy greeted by banner begging him to dress for the occassion did he hulk hogan couldnt stand gavin newsom exwife say about govs online trolling rightwing reporter brian glenn praise zelenskys attire after shaming him in february texas democrat return ending twoweek walkout over redistricting trump say hell talk to putin after meeting with zelenskyy and european leader at the white house reporter brian glenn brag about dc being safe after walking street with girlfriend marjorie taylor greene zelenskyy trump express hope for trilateral talk with putin to bring end to russiaukraine war oval office meeting again put rocky tense trumpzelenskyy relationship to the test armstrong williams how george washington would end russia war on ukraine staff commentary la habra man charged in fatal whittier shooting of driver texas democrat head home after week in illinois mpo holding first public planning session tuesday at mountain line photo show aftermath of a deadly russian attack in ukraine and the funeral of a fallen soldier now thats a selfown eric swalwell grocery shopping go viral for all the wrong hilarious reason texas democrat return to capitol for republican u house map overhaul a california counter texas democrat return to capitol for republican u house map overhaul a california counter letter time for democrat to play to win volodymyr zelensky wear suit this time for white house meeting with trump tarrant county commissioner again propose cutting voting site ahead of an election this is exactly what technology should be doing man shot injured in parking lot of huntington beach store letter no place for threat violence in politics trump aim to set putinzelenskiy meeting after monday summit obbba creates a cohesive framework for u hydrogen growth top trump lackey demand dc resident kiss the ground for occupation republican look to make a uturn on electric vehicle for the postal service republican look to make a uturn on electric vehicle for the postal service"
$ws.Range("J40").Value = "poriferous llc achieves coveted eu mdr regulatory certification porous polyethylene implant warehouse data reveals shipper managed inventory tariff nvidia earnings face high bar datacenter demand china deal newsmax pay million settle dominion suit election fraud claim bmws next big suv gunning dog electric jules witcover political reporter columnist dy samsung overtakes apple smartphone market update hint major shift brent wti edge higher zelensky arrives washington trump trending coin oz presale token gain early momentum best altcoins buy september remittix list ada pepe bonk cold wallet offer referral bonus bonk burn supply ripple eye crypto project set explode crypto executive seen accumulating million token new meme coin challenging shiba inu shib walmart selling queen mattress topper thats soft cool touch best electronic signature apps judge denies bail north andover offduty officer shot defense file another appeal whale bet big new payment token replace xrp stellar carvanas breakout fuel interest consumer ecommerce etf colorado rabbit growing tentacle horn maluma sign wme area sky limit dow jones today stock tick lower major index hover record high consecutive gain investopedia fake labubu doll pose threat kid safety watchdog cbs news court federal claim small business landmark taa v baa protest employer guidance new tax deduction tip overtime pay terra balcanica intersects gold mineralization start drill program brezani target bosnia business summarybrief pm conservative network newsmax agrees pay defamation case bogus election claim invesco qqq qqq share sold lee financial co bank stock eye august terawulf stock jump google boost stake data center operator bse gainer watch uncy investor opportunity lead unicycive therapeutic security fraud lawsuit schall law firm cornerstone select advisor llc reduces stock pfizer pfe southeast toyota finance ranked highest overall dealer satisfaction jd power forget ai robotic process automation healthcare hottest technology warehouse teamster unfi lucrative contract promising retail stock follow august alphabet goog share sold wellington management llp promising lithium stock eye august insurance stock watch today august hydrogen stock consider august pwc kunai acquisition bolster tech enablement bank cfo bull case nvidia got stronger according analyst curved samsung screen lowest price today soho house taken private billion deal new york expect fiery outcome fed jackson hole conference financial advisor network take cocacola company ko analyst rating campbell mesa laboratory executes previouslyannounced strategic financing plan ibm option trading deep dive market sentiment forecasting future analyst projection cloudflare closer look walmarts option market dynamic flex meal deal faculty staff graduate professional student paypal holdingss option frenzy glimpse expert outlook jm smucker analyst assessing microsoft insight financial analyst ethereum trader turned million month analyst verdict fedex eye expert ai take apple edge foxconns pivot spark etf shift prudential financial stand analyst mesa laboratory executes previouslyannounced strategic financing plan wired roundup gpt flopped wired market bullish bearish royal gold willis tower watson recent short interest look pge price earnings market feeling lumentum holding ozempic creator novo nordisk sell weight loss drug cheaper new partnership novo nordisk offer ozempic per month eligible cashpaying customer reuters texas shut data center protect energy grid blockdag reach analyst review presale dogecoin news today dogecoin price decline technical indicator suggest potential rebound applied digitals stock apld jump plan billion ai campus tipranks worker job hugging stagnant labor market growing resentment bail soon next great resignation come yahoocom texas attorney general accuses meta characterai misleading kid mental health claim techcrunch worker job hugging stagnant labor market growing resentment bail soon next great resignation come fortune conservative network newsmax agrees pay defamation case bogus election claim skechers named title sponsor world champion cup shriners childrens announced charity partner feather sound country club clearwater florida host tournament whats propanc share monday mastercard worldpay team streamline uae money movement mcdonalds corporation mcd share sold community bank trust waco texas granite harbor advisor boost holding mcdonalds corporation mcd walmart everything head school backpack laptop madd join force law enforcement prevent impaired driving ahead labor day saturation saturday initiative weekend ibd qa stock list monday aug spotlight fed powell bowman waller rate dogecoin rally spark buzz trader sub coin upside fairscale capital llc purchase share crowdstrike crwd lee financial co sell share danaher corporation dhr berger montague pc investigates security claim kindercare learning company nyse klc amazon selling bluetooth speaker perfect setting inventhelp inventor develops new training aid basketball player chk data science recruiter find ai analytics ml expert gsg advisor llc purchase share invesco qqq qqq gsg advisor llc stock holding mastercard incorporated weekly bankruptcy alert august ending august jackson hole loom powell speech test september ratecut newsmax pay million settle dominion defamation lawsuit intel stock drop report equity stake via chip act barrons liquid edge ai platform leap expands support laptop bestinclass performance amd ryzentm ryzen aitm processor business summarybrief pm business summarybrief pm texas democrat froma harrop parent dont vaccinate child guilty abuse handwritten proof holocaust theft compel congress act news summary pm news summary pm dc sue doj stop immigration requirement victim fund trump zelenskyy meeting end road support securing peace deal editorial roundup united eric adam roll nypd expansion trump dc crackdown factor trending summarybrief pm texas lawmaker illinois heading following twoweek walkout remap fight nfl team refuse facing backlash hiring male cheerleader stephen miller trump dc crackdown target graffiti news summary pm conservative network newsmax agrees pay defamation case bogus election claim mark hamill wife clever thought stopped fleeing trump trump zelensky put smile oval office vancouver mayor warns delay bridge replacement cost per day trump zelensky put smile oval office hamas accepts arab ceasefire proposal gaza palestinian death toll pas texas redistricting fight shift democrat austin new strategy mississippi becomes fourth send national guard troop dc expanding federal crackdown mississippi becomes fourth send national guard troop dc expanding federal crackdown mississippi becomes fourth send national guard troop dc expanding federal crackdown mississippi becomes fourth send national guard troop dc expanding federal crackdown mississippi becomes fourth send national guard troop dc expanding federal crackdown mississippi becomes fourth send national guard troop dc expanding federal crackdown mississippi becomes fourth send national guard troop dc expanding federal crackdown mississippi becomes fourth send national guard troop dc expanding federal crackdown mississippi becomes fourth send national guard troop dc expanding federal crackdown mississippi becomes fourth send national guard troop dc expanding federal crackdown putin advice trump launch assault mailin ballot voting machine eo seek expand access alternative asset investment retirement plan doj demand sensitive illinois voter registration data mike mcgrew alternative mindset patriotic independent commentary brian stelter doesnt making msnbcs rebranding funnier excited thread attorney general william barr didnt trump implicated epstein file attorney general william barr didnt trump implicated epstein file gavin newsom left fox news speechless unexpectedly hilarious woke comeback line defense trump abuse power subsidence isnt environmental crisis slowed latest zelenskyy signal openness threeway meeting trump putin doj wont appeal judge order concealed carry ban post office latest zelenskyy signal openness threeway meeting trump putin letter political treachery versus personal charity redistricting fight texas democrat california start work redistricting fight texas democrat california start work redistricting fight texas democrat california start work redistricting fight texas democrat california start work redistricting fight texas democrat california start work redistricting fight texas democrat california start work redistricting fight texas democrat california start work redistricting fight texas democrat california start work redistricting fight texas democrat california start work bnpl provider extension nydfs request information deadline wyoming teen educate youth human trafficking social medium congress must injustice facing combatinjured veteran blessing disaster shape alabama c galaxy iconic aircraft texas attorney general accuses meta characterai misleading kid mental health claim texas ag accuses meta characterai misleading kid mental health claim european join trump zelenskyy white house conservative network newsmax agrees pay defamation case bogus election claim redistricting fight texas democrat california start work redistricting fight texas democrat california start work twoweek standoff end quorum reached texas house mississippi becomes fourth send national guard troop dc expanding federal crackdown letter trumpputin photo op accomplished nothing judge weighs detainee legal alligator alcatraz florida everglades judge weighs detainee legal alligator alcatraz florida everglades roblox locking sexual content access adult location lawsuit trump european security guarantee ukraine help end war russia good crisis dr nirav shah promising candidate anything steve collins wear suit bro zelenskyy greeted banner begging dress occassion hulk hogan couldnt stand gavin newsom exwife govs online trolling rightwing reporter brian glenn praise zelenskys attire shaming february texas democrat ending twoweek walkout redistricting trump hell putin meeting zelenskyy european white house reporter brian glenn brag dc safe walking street girlfriend marjorie taylor greene zelenskyy trump express trilateral putin end russiaukraine war oval office meeting put rocky tense trumpzelenskyy relationship test armstrong williams george washington end russia war ukraine staff commentary la habra man charged fatal whittier shooting driver texas democrat head illinois mpo holding public planning session tuesday mountain line photo aftermath deadly russian attack ukraine funeral fallen soldier thats selfown eric swalwell grocery shopping viral hilarious texas democrat capitol republican house map overhaul california counter texas democrat capitol republican house map overhaul california counter letter democrat play volodymyr zelensky wear suit white house meeting trump tarrant county commissioner propose cutting voting site ahead election exactly technology man shot injured parking lot huntington beach store letter place threat violence politics trump aim set putinzelenskiy meeting monday summit obbba creates cohesive framework hydrogen growth trump lackey demand dc resident kiss ground occupation republican look uturn electric vehicle postal service republican look uturn electric vehicle postal service"

Write-Output "Row 40 written"
